# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.557.73"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.903.33"
$ws.Range("E3").Value = "  -2.86%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.13%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").Value = "2.911.06"
$ws.Range("E9").Value = "  -2.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.359"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.22%  "
$ws.Range("D13").Value = "3.403.58"
$ws.Range("E13").Value = "  -3.21%  "
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "60.519.70"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.51%  "
$ws.Range("D17").Value = "2.908.42"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000141"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "363.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "3.026.10"
$ws.Range("E25").Value = "  -3.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.179"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.61%  "
$ws.Range("D30").Value = "0.0₃0857"
$ws.Range("E30").Value = "  -8.40%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.40%  "
$ws.Range("E38").Value = "  -6.52%  "
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "2.327.72"
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.55%  "
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0572"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0234"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0931"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.08%  "
